$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (values like "1.00", "7.23", etc.)
# so they stay text, matching the rest of the Price column.
$textCells = @("D4", "D5", "D6", "D14", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D33", "D34", "D35", "D37", "D39", "D40", "D43", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (prices and hourly % changes),
# including the dogwifhat / VeChain row swap (rows 49-50).
$ws.Range("D2").Value = "61.039.92"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.386.10"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "572.31"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "141.38"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").Value = "3.966.22"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "27.89"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "3.387.61"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "61.132.71"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").Value = "8.97"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "384.79"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "76.11"
$ws.Range("E22").Value = "  +3.51%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "0.0000115"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "0.185"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("D33").Value = "23.32"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").Value = "6.96"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "165.69"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").Value = "3.420.14"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "5.00"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("D39").Value = "0.0767"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").Value = "26.62"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "4.36"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "2.461.03"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "22.88"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "6.65"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0263"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.13"
$ws.Range("E50").Value = "  +10.13%  "
$ws.Range("D51").Value = "0.206"
$ws.Range("E51").Value = "  -1.87%  "
